$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# 1. Rename task in row 14 ("tela" -> "senha")
$ws.Range("B14").Value2 = "Criar tela de alteração de senha"

# 2. Mark "Criar tela livro Caixa" (row 12) as "Stand by"
$ws.Range("D12").Value2 = "Stand by"

# 3. Mark the three "contas a pagar e receber" rows as done ("ok")
$ws.Range("D17").Value2 = "ok"
$ws.Range("D18").Value2 = "ok"
$ws.Range("D19").Value2 = "ok"

# 4. Insert a new row before the old row 21, pushing the "Banco de dados" block
#    (old rows 21-59) down to rows 22-60, leaving row 20 free for the new task.
$ws.Rows.Item(21).Insert()

# 5. Populate new row 20 with the "Login" front-end task
$ws.Range("B20").Value2 = "Login"
$ws.Range("C20").Value2 = 0
$ws.Range("D20").Value2 = "ok"

# Copy number formats from neighbouring cells so the new row matches the
# sheet's existing look (time format on C, date format on D), reusing the
# existing style entries rather than creating new ones.
$ws.Range("C19").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# 6. Update the active selection to B4, as in the saved workbook
$ws.Range("B4").Select()
